# Scheduled-runner style refresh of the market-price-derived profit
# columns (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ, columns H-N) across several leve rows in
# several job sheets of the Mandragora_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 499
$ws.Range("I31").Value = 499
$ws.Range("K31").Value = 1497
$ws.Range("M31").Value = -1267
$ws.Range("H40").Value = 2851
$ws.Range("J40").Value = 2461.8572
$ws.Range("L40").Value = 2461.8572
$ws.Range("N40").Value = -2811.8572
$ws.Range("H43").Value = 895
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 661.1111
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 661.1111
$ws.Range("M43").Value = -2931
$ws.Range("N43").Value = -799.1111
$ws.Range("H46").Value = 33334914
$ws.Range("J46").Value = 1975
$ws.Range("L46").Value = 5925
$ws.Range("N46").Value = -6163
$ws.Range("H60").Value = 33334914
$ws.Range("J60").Value = 1975
$ws.Range("L60").Value = 5925
$ws.Range("N60").Value = -6893
$ws.Range("H64").Value = 3920.2856
$ws.Range("I64").Value = 3492.8572
$ws.Range("J64").Value = 4775.143
$ws.Range("K64").Value = 3492.8572
$ws.Range("L64").Value = 4775.143
$ws.Range("M64").Value = -3244.8572
$ws.Range("N64").Value = -5271.143
$ws.Range("H67").Value = 3920.2856
$ws.Range("I67").Value = 3492.8572
$ws.Range("J67").Value = 4775.143
$ws.Range("K67").Value = 3492.8572
$ws.Range("L67").Value = 4775.143
$ws.Range("M67").Value = -2634.8572
$ws.Range("N67").Value = -6491.143
$ws.Range("H76").Value = 3164.5
$ws.Range("I76").Value = 2829
$ws.Range("K76").Value = 2829
$ws.Range("M76").Value = -2514
$ws.Range("H79").Value = 3164.5
$ws.Range("I79").Value = 2829
$ws.Range("K79").Value = 2829
$ws.Range("M79").Value = -1737
$ws.Range("H138").Value = 4030.1538
$ws.Range("I138").Value = 1480.4736
$ws.Range("J138").Value = 5498.1514
$ws.Range("K138").Value = 4441.4208
$ws.Range("L138").Value = 16494.4542
$ws.Range("M138").Value = 698.5792000000001
$ws.Range("N138").Value = -26774.4542

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 3700
$ws.Range("I35").Value = 3700
$ws.Range("K35").Value = 3700
$ws.Range("M35").Value = -3294
$ws.Range("H45").Value = 3004.1785
$ws.Range("I45").Value = 1804.1666
$ws.Range("J45").Value = 3904.1875
$ws.Range("K45").Value = 1804.1666
$ws.Range("L45").Value = 3904.1875
$ws.Range("M45").Value = -1427.1666
$ws.Range("N45").Value = -4658.1875
$ws.Range("H61").Value = 2769.3845
$ws.Range("I61").Value = 1501.7059
$ws.Range("K61").Value = 1501.7059
$ws.Range("M61").Value = -1289.7059
$ws.Range("H136").Value = 2769.3845
$ws.Range("I136").Value = 1501.7059
$ws.Range("K136").Value = 4505.1177
$ws.Range("M136").Value = -1955.1177

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 751.41174
$ws.Range("I22").Value = 947.6667
$ws.Range("J22").Value = 280.4
$ws.Range("K22").Value = 947.6667
$ws.Range("L22").Value = 280.4
$ws.Range("M22").Value = -774.6667
$ws.Range("N22").Value = -626.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 7096
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 7620
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 7620
$ws.Range("M4").Value = -4888
$ws.Range("N4").Value = -7844
$ws.Range("H132").Value = 5886256
$ws.Range("I132").Value = 9093578
$ws.Range("J132").Value = 6166.3335
$ws.Range("K132").Value = 27280734
$ws.Range("L132").Value = 18499.0005
$ws.Range("M132").Value = -27278204
$ws.Range("N132").Value = -23559.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2879348
$ws.Range("I4").Value = 3636929.2
$ws.Range("J4").Value = 539.6
$ws.Range("K4").Value = 10910787.6
$ws.Range("L4").Value = 1618.8
$ws.Range("M4").Value = -10910675.6
$ws.Range("N4").Value = -1842.8
$ws.Range("H131").Value = 735926.5
$ws.Range("I131").Value = 515.6429000000001
$ws.Range("J131").Value = 4167844
$ws.Range("K131").Value = 1546.9287
$ws.Range("L131").Value = 12503532
$ws.Range("M131").Value = 3493.0713
$ws.Range("N131").Value = -12513612

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 37502.5
$ws.Range("H29").Value = 10003000
$ws.Range("I29").Value = 10003000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 10003000
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -10002710
$ws.Range("N29").ClearContents()
$ws.Range("H70").Value = 5229.8125
$ws.Range("I70").Value = 5211.8
$ws.Range("J70").Value = 5500
$ws.Range("K70").Value = 5211.8
$ws.Range("L70").Value = 5500
$ws.Range("M70").Value = -4941.8
$ws.Range("N70").Value = -6040
$ws.Range("H73").Value = 5229.8125
$ws.Range("I73").Value = 5211.8
$ws.Range("J73").Value = 5500
$ws.Range("K73").Value = 5211.8
$ws.Range("L73").Value = 5500
$ws.Range("M73").Value = -4275.8
$ws.Range("N73").Value = -7372
$ws.Range("H80").Value = 3147.5
$ws.Range("I80").Value = 2981
$ws.Range("J80").Value = 3980
$ws.Range("K80").Value = 2981
$ws.Range("L80").Value = 3980
$ws.Range("M80").Value = -1983
$ws.Range("N80").Value = -5976
$ws.Range("H83").Value = 3147.5
$ws.Range("I83").Value = 2981
$ws.Range("J83").Value = 3980
$ws.Range("K83").Value = 14905
$ws.Range("L83").Value = 19900
$ws.Range("M83").Value = -9913
$ws.Range("N83").Value = -29884
$ws.Range("H126").Value = 4145.7
$ws.Range("I126").Value = 2719.4546
$ws.Range("J126").Value = 5888.8887
$ws.Range("K126").Value = 8158.3638
$ws.Range("L126").Value = 17666.6661
$ws.Range("M126").Value = -5688.3638
$ws.Range("N126").Value = -22606.6661
$ws.Range("H132").Value = 1661.3636
$ws.Range("I132").Value = 1438.2941
$ws.Range("J132").Value = 2419.8
$ws.Range("K132").Value = 4314.8823
$ws.Range("L132").Value = 7259.400000000001
$ws.Range("M132").Value = -1784.8823
$ws.Range("N132").Value = -12319.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H82").Value = 2867.1667
$ws.Range("I82").Value = 4400
$ws.Range("J82").Value = 2560.6
$ws.Range("K82").Value = 4400
$ws.Range("L82").Value = 2560.6
$ws.Range("M82").Value = -4039
$ws.Range("N82").Value = -3282.6
$ws.Range("H85").Value = 2867.1667
$ws.Range("I85").Value = 4400
$ws.Range("J85").Value = 2560.6
$ws.Range("K85").Value = 4400
$ws.Range("L85").Value = 2560.6
$ws.Range("M85").Value = -3152
$ws.Range("N85").Value = -5056.6
$ws.Range("H136").Value = 16130364
$ws.Range("I136").Value = 26316922
$ws.Range("K136").Value = 78950766
$ws.Range("M136").Value = -78948216
